$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Remove the "sex" column (column E) from the test_file sheet, shifting
# customer_type left into its place.
$ws1.Columns.Item(5).Delete()

# Duplicate the two sample data rows on Sheet2 twice more (rows 4-5, 6-7).
$ws2.Range("A2:C3").Copy()
$ws2.Range("A4").PasteSpecial()
$ws2.Range("A2:C3").Copy()
$ws2.Range("A6").PasteSpecial()

# Update selections: Sheet2's selection moves to F16 and is no longer the
# active tab; the test_file sheet becomes active with I8 selected.
$ws2.Range("F16").Select()
$ws1.Activate()
$ws1.Range("I8").Select()
